# Update Users sheet: LastLogin timestamp for smfasihaly
$wb = $excel.ActiveWorkbook

$wsUsers = $wb.Worksheets.Item("Users")
$wsUsers.Range("C2").Value = "2024-06-24 08:48:02"

# JustFlipped sheet: remove rows 5-18 (old extra vocab rows), keep rows 1-4,
# and update rows 2-4 with the new words / translations / direction.
$wsFlipped = $wb.Worksheets.Item("JustFlipped")
$wsFlipped.Rows.Item(5).Resize(14).Delete()

$wsFlipped.Range("A2").Value = "cadere"
$wsFlipped.Range("B2").Value = "to fall"
$wsFlipped.Range("D2").Value = "Italian to English"

$wsFlipped.Range("A3").Value = "scendere"
$wsFlipped.Range("B3").Value = "to go down"
$wsFlipped.Range("D3").Value = "Italian to English"

$wsFlipped.Range("A4").Value = "tenere"
$wsFlipped.Range("B4").Value = "to keep"
$wsFlipped.Range("D4").Value = "Italian to English"

# Failure sheet: append 4 new rows (11-14) of vocab entries
$wsFailure = $wb.Worksheets.Item("Failure")

$wsFailure.Range("A11").Value = "ridere"
$wsFailure.Range("B11").Value = "to laugh"
$wsFailure.Range("C11").Value = "smfasihaly"
$wsFailure.Range("D11").Value = "Italian to English"

$wsFailure.Range("A12").Value = "chiudere"
$wsFailure.Range("B12").Value = "to close"
$wsFailure.Range("C12").Value = "smfasihaly"
$wsFailure.Range("D12").Value = "Italian to English"

$wsFailure.Range("A13").Value = "imparare"
$wsFailure.Range("B13").Value = "to learn"
$wsFailure.Range("C13").Value = "smfasihaly"
$wsFailure.Range("D13").Value = "Italian to English"

$wsFailure.Range("A14").Value = "incontrare"
$wsFailure.Range("B14").Value = "to encounter"
$wsFailure.Range("C14").Value = "smfasihaly"
$wsFailure.Range("D14").Value = "Italian to English"
